$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2026-02-18 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-19 Thursday", 2)

# Update the division-problem answers in the first table. Cells are
# addressed by (row, column) rather than by text search because several
# answer strings (e.g. "70÷3=23, 1") appear more than once in the table
# and need to change to different new values depending on position.
$tbl = $d.Tables.Item(1)

$newValues = @{
    1 = @("60÷2=30, 0", "91÷7=13, 0", "99÷8=12, 3", "12÷8=1, 4", "95÷2=47, 1")
    5 = @("17÷5=3, 2", "56÷4=14, 0", "84÷2=42, 0", "58÷8=7, 2", "79÷2=39, 1")
    9 = @("65÷9=7, 2", "41÷4=10, 1", "17÷5=3, 2", "40÷5=8, 0", "51÷2=25, 1")
    13 = @("47÷7=6, 5", "85÷4=21, 1", "29÷2=14, 1", "80÷5=16, 0", "97÷7=13, 6")
    17 = @("14÷2=7, 0", "78÷7=11, 1", "35÷8=4, 3", "85÷7=12, 1", "15÷2=7, 1")
}

foreach ($rowIndex in $newValues.Keys) {
    $rowValues = $newValues[$rowIndex]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $cell = $tbl.Cell($rowIndex, $col)
        $cell.Range.Text = $rowValues[$col - 1]
    }
}
